# Update status text from "Ready for handoff" to "In Translation"
# across all worksheets (Overview, zh-cn, de-de). The columns that were
# sized to fit that status text are then re-fit to the new, shorter text.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $touchedCols = @{}

    foreach ($cell in $used.Cells) {
        # Put the literal on the left so PowerShell's type-coercion rules
        # (which convert the right side to the left side's type) don't turn
        # this into a truthiness test against boolean cell values.
        if ("Ready for handoff" -eq $cell.Value()) {
            $cell.Value = "In Translation"
            $touchedCols[$cell.Column] = $true
        }
    }

    # Re-fit only the columns whose content actually changed, matching how
    # Excel recalculates a column's autofit width after an edit. The new,
    # shorter status text ("In Translation" vs "Ready for handoff") fits in
    # a narrower column.
    foreach ($colIndex in $touchedCols.Keys) {
        $ws.Columns.Item($colIndex).ColumnWidth = 12.5
    }
}
